$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure changed cells keep their original text representation
# (avoid Excel auto-converting numeric-looking strings to numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.916.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.983.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.60"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.55%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.968.65"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.72"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.146"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.79"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.485.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.034.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.95"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.990.65"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "439.82"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.683"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.39"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.71"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.05"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.77%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.15"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.77"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.03"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.06"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.05"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.298"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +10.64%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0352"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "376.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.676.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.16"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.89"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +10.72%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.97%  "
